# Edit: split the single "HYDROGEN" storage commodity row (row 14) into a
# HYDROGEN_IN (row 14) / HYDROGEN_OUT (row 15) pair on the "data" sheet, and
# repoint the output-commodity formula (G7) at the new row. Also stamp the
# NRG set flag on H20 and move the active selection to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# 1) Row 14: commodity name HYDROGEN -> HYDROGEN_IN
$ws.Range("C14").Value = "HYDROGEN_IN"

# 2) Duplicate row 14's formatting into row 15 (same style s="17" on B:H)
$ws.Range("B14:H14").Copy()
$ws.Range("B15:H15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill in row 15 values for the new HYDROGEN_OUT commodity membership row
$ws.Range("B15").Value = "NRG"
$ws.Range("C15").Value = "HYDROGEN_OUT"
$ws.Range("D15").Value = "wodor"
$ws.Range("E15").Value = "PJ"
$ws.Range("G15").Value = "DAYNITE"
$ws.Range("H15").Value = "ANNUAL"

# 4) G7 (output commodity) now looks up the new row instead of reusing row 14
$ws.Range("G7").Formula = "=C15"

# 5) H20 gains the NRG set flag
$ws.Range("H20").Value = "NRG"

# 6) Move the active selection to G8 (matches the saved selection in the file)
$ws.Range("G8").Select()
